$d = $word.ActiveDocument

# Fix #366 - User content is lost after two generation without edition.
# Word "simplifies" a field that has never been updated back down to a
# <w:fldSimple> on save; m2doc's user-content markers (m:usercontent /
# m:endusercontent) must instead be stored as a "complex" field
# (separate begin/instrText/separate/end runs) so the content between
# them survives a second round of generation. Walk every paragraph,
# find the ones still holding a simple field, and rewrite them using
# the begin/instrText/separate/end run sequence, preserving the
# paragraph's own formatting/attributes and the field instruction text.

$pkgHeader = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>'
$pkgFooter = '</pkg:xmlData></pkg:part></pkg:package>'
$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $xml = $p.Range.WordOpenXML

    if ($xml -match '<w:p\b([^>]*)>\s*<w:fldSimple w:instr="([^"]*)"\s*/?>\s*(?:</w:fldSimple>)?\s*</w:p>') {
        $attrs = $matches[1]
        $instr = $matches[2]

        # Drop the w14:paraId/textId attributes the OOXML round-trip adds;
        # they are not present on the original paragraph mark.
        $attrs = $attrs -replace '\s*w14:paraId="[^"]*"', ''
        $attrs = $attrs -replace '\s*w14:textId="[^"]*"', ''

        $newParagraph = '<w:p' + $attrs + '>' `
            + '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' `
            + '<w:r><w:instrText>' + $instr + '</w:instrText></w:r>' `
            + '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' `
            + '<w:r><w:fldChar w:fldCharType="end"/></w:r>' `
            + '</w:p>'

        $body = '<w:document ' + $wordNs + '><w:body>' + $newParagraph + '<w:sectPr/></w:body></w:document>'

        $p.Range.InsertXML($pkgHeader + $body + $pkgFooter)
    }
}
